$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 17-30: Fecha (D), Volumen (M), Precio minimo (N), Precio maximo (O), Precio promedio ponderado (P), Precio $/Kg (S)
$ws.Range("D17").Value = 44781
$ws.Range("M17").Value = 160
$ws.Range("N17").Value = 23000
$ws.Range("O17").Value = 24000
$ws.Range("P17").Value = 23500
$ws.Range("S17").Value = 1175

$ws.Range("D18").Value = 44778
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = 23000
$ws.Range("O18").Value = 24000
$ws.Range("P18").Value = 23500
$ws.Range("S18").Value = 1175

$ws.Range("D19").Value = 44435
$ws.Range("M19").Value = 260
$ws.Range("N19").Value = 20000
$ws.Range("O19").Value = 22000
$ws.Range("P19").Value = 21115
$ws.Range("S19").Value = 1056

$ws.Range("D20").Value = 44343
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 19500
$ws.Range("O20").Value = 20000
$ws.Range("P20").Value = 19750
$ws.Range("S20").Value = 988

$ws.Range("D21").Value = 44364
$ws.Range("M21").Value = 140
$ws.Range("N21").Value = 20000
$ws.Range("O21").Value = 21000
$ws.Range("P21").Value = 20500
$ws.Range("S21").Value = 1025

$ws.Range("D22").Value = 44428
$ws.Range("M22").Value = 100
$ws.Range("N22").Value = 20000
$ws.Range("O22").Value = 21000
$ws.Range("P22").Value = 20500
$ws.Range("S22").Value = 1025

$ws.Range("D23").Value = 44333
$ws.Range("M23").Value = 100
$ws.Range("N23").Value = 19500
$ws.Range("O23").Value = 20000
$ws.Range("P23").Value = 19750
$ws.Range("S23").Value = 988

$ws.Range("D24").Value = 44431
$ws.Range("M24").Value = 160
$ws.Range("N24").Value = 21000
$ws.Range("O24").Value = 22000
$ws.Range("P24").Value = 21500
$ws.Range("S24").Value = 1075

$ws.Range("D25").Value = 44420
$ws.Range("M25").Value = 160
$ws.Range("N25").Value = 20000
$ws.Range("O25").Value = 21000
$ws.Range("P25").Value = 20500
$ws.Range("S25").Value = 1025

$ws.Range("D26").Value = 44365
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 20000
$ws.Range("O26").Value = 21000
$ws.Range("P26").Value = 20500
$ws.Range("S26").Value = 1025

$ws.Range("D27").Value = 44417
$ws.Range("M27").Value = 160
$ws.Range("N27").Value = 20000
$ws.Range("O27").Value = 21000
$ws.Range("P27").Value = 20500
$ws.Range("S27").Value = 1025

$ws.Range("D28").Value = 44427
$ws.Range("M28").Value = 200
$ws.Range("N28").Value = 20000
$ws.Range("O28").Value = 21000
$ws.Range("P28").Value = 20500
$ws.Range("S28").Value = 1025

$ws.Range("D29").Value = 44441
$ws.Range("M29").Value = 160
$ws.Range("N29").Value = 20000
$ws.Range("O29").Value = 21000
$ws.Range("P29").Value = 20500
$ws.Range("S29").Value = 1025

$ws.Range("D30").Value = 44434
$ws.Range("M30").Value = 100
$ws.Range("N30").Value = 20000
$ws.Range("O30").Value = 21000
$ws.Range("P30").Value = 20500
$ws.Range("S30").Value = 1025

# Append new rows 31-33 with full data, replicating the constant columns from the existing dataset
$ws.Range("A31").Value = 8
$ws.Range("B31").Value = "Terminal La Palmera de La Serena"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44776
$ws.Range("D31").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100108
$ws.Range("H31").Value = "Tropicales y subtropicales"
$ws.Range("I31").Value = 100108007
$ws.Range("J31").Value = "Coco"
$ws.Range("K31").Value = "Sin especificar"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 160
$ws.Range("N31").Value = 23000
$ws.Range("O31").Value = 24000
$ws.Range("P31").Value = 23500
$ws.Range("Q31").Value = "$/malla 20 unidades"
$ws.Range("R31").Value = "Perú"
$ws.Range("S31").Value = 1175
$ws.Range("T31").Value = 20

$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = 44301
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100108
$ws.Range("H32").Value = "Tropicales y subtropicales"
$ws.Range("I32").Value = 100108007
$ws.Range("J32").Value = "Coco"
$ws.Range("K32").Value = "Sin especificar"
$ws.Range("L32").Value = "Primera"
$ws.Range("M32").Value = 100
$ws.Range("N32").Value = 18000
$ws.Range("O32").Value = 19000
$ws.Range("P32").Value = 18500
$ws.Range("Q32").Value = "$/malla 20 unidades"
$ws.Range("R32").Value = "Perú"
$ws.Range("S32").Value = 925
$ws.Range("T32").Value = 20

$ws.Range("A33").Value = 8
$ws.Range("B33").Value = "Terminal La Palmera de La Serena"
$ws.Range("C33").Value = "Coquimbo"
$ws.Range("D33").Value = 44336
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E33").Value = 4
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100108
$ws.Range("H33").Value = "Tropicales y subtropicales"
$ws.Range("I33").Value = 100108007
$ws.Range("J33").Value = "Coco"
$ws.Range("K33").Value = "Sin especificar"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 100
$ws.Range("N33").Value = 19500
$ws.Range("O33").Value = 20000
$ws.Range("P33").Value = 19750
$ws.Range("Q33").Value = "$/malla 20 unidades"
$ws.Range("R33").Value = "Perú"
$ws.Range("S33").Value = 988
$ws.Range("T33").Value = 20
